$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 113 (hunk @ diff line 6316)
$ws.Range("H113").Value = 4907.7856
$ws.Range("I113").Value = 4609.909
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 4609.909
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -1355.909
$ws.Range("N113").Value = -12508

# row 137 (hunk @ diff line 7522)
$ws.Range("H137").Value = 3354.8
$ws.Range("J137").Value = 2500.7144
$ws.Range("L137").Value = 7502.1432
$ws.Range("N137").Value = -12602.1432

# row 138 (hunk @ diff line 7574)
$ws.Range("H138").Value = 4158.8657
$ws.Range("I138").Value = 2877.7144
$ws.Range("J138").Value = 4422.6323
$ws.Range("K138").Value = 8633.143199999999
$ws.Range("L138").Value = 13267.8969
$ws.Range("M138").Value = -3493.143199999999
$ws.Range("N138").Value = -23547.8969

$ws = $wb.Worksheets.Item("ARM")
# row 32 (hunk @ diff line 9358)
$ws.Range("H32").Value = 6371.2466
$ws.Range("I32").Value = 4593.067
$ws.Range("J32").Value = 14578.23
$ws.Range("K32").Value = 4593.067
$ws.Range("L32").Value = 14578.23
$ws.Range("M32").Value = -4306.067
$ws.Range("N32").Value = -15152.23

# row 37 (hunk @ diff line 9615)
$ws.Range("H37").Value = 20585.5
$ws.Range("I37").Value = 7000
$ws.Range("J37").Value = 23302.6
$ws.Range("K37").Value = 7000
$ws.Range("L37").Value = 23302.6
$ws.Range("M37").Value = -6727
$ws.Range("N37").Value = -23848.6

# row 44 (hunk @ diff line 9967)
$ws.Range("H44").Value = 22000
$ws.Range("J44").Value = 22000
$ws.Range("L44").Value = 22000
$ws.Range("N44").Value = -22976

# row 74 (hunk @ diff line 11449)
$ws.Range("H74").Value = 784.8444
$ws.Range("I74").Value = 735.13336
$ws.Range("J74").Value = 884.26666
$ws.Range("K74").Value = 735.13336
$ws.Range("L74").Value = 884.26666
$ws.Range("M74").Value = 138.86664
$ws.Range("N74").Value = -2632.26666

# row 77 (hunk @ diff line 11596)
$ws.Range("H77").Value = 784.8444
$ws.Range("I77").Value = 735.13336
$ws.Range("J77").Value = 884.26666
$ws.Range("K77").Value = 3675.6668
$ws.Range("L77").Value = 4421.3333
$ws.Range("M77").Value = 692.3332
$ws.Range("N77").Value = -13157.3333

# row 132 (hunk @ diff line 14288)
$ws.Range("H132").Value = 14287848
$ws.Range("I132").Value = 20409334
$ws.Range("K132").Value = 61228002
$ws.Range("M132").Value = -61225472

# row 137 (hunk @ diff line 14533)
$ws.Range("H137").Value = 29608.889
$ws.Range("J137").Value = 29608.889
$ws.Range("L137").Value = 29608.889
$ws.Range("N137").Value = -39808.889

$ws = $wb.Worksheets.Item("BSM")
# row 20 (hunk @ diff line 15769)
$ws.Range("H20").Value = 1538.7778
$ws.Range("I20").Value = 910.2353000000001
$ws.Range("J20").Value = 2607.3
$ws.Range("K20").Value = 910.2353000000001
$ws.Range("L20").Value = 2607.3
$ws.Range("M20").Value = -663.2353000000001
$ws.Range("N20").Value = -3101.3

$ws = $wb.Worksheets.Item("CRP")
# row 31 (hunk @ diff line 23283)
$ws.Range("H31").Value = 4672.3057
$ws.Range("I31").Value = 3265.6086
$ws.Range("J31").Value = 7161.077
$ws.Range("K31").Value = 3265.6086
$ws.Range("L31").Value = 7161.077
$ws.Range("M31").Value = -2970.6086
$ws.Range("N31").Value = -7751.077

# row 34 (hunk @ diff line 23430)
$ws.Range("H34").Value = 4672.3057
$ws.Range("I34").Value = 3265.6086
$ws.Range("J34").Value = 7161.077
$ws.Range("K34").Value = 3265.6086
$ws.Range("L34").Value = 7161.077
$ws.Range("M34").Value = -3063.6086
$ws.Range("N34").Value = -7565.077

# row 58 (hunk @ diff line 24609)
$ws.Range("H58").Value = 9261668
$ws.Range("J58").Value = 27782106
$ws.Range("L58").Value = 27782106
$ws.Range("N58").Value = -27782512

# row 95 (hunk @ diff line 26425)
$ws.Range("H95").Value = 15764
$ws.Range("J95").Value = 15764
$ws.Range("L95").Value = 15764
$ws.Range("N95").Value = -21256

# row 96 (hunk @ diff line 26474)
$ws.Range("H96").Value = 20854
$ws.Range("J96").Value = 20854
$ws.Range("L96").Value = 20854
$ws.Range("N96").Value = -26346

# row 132 (hunk @ diff line 28238)
$ws.Range("H132").Value = 2059.5
$ws.Range("I132").Value = 1565.027
$ws.Range("J132").Value = 3722.7273
$ws.Range("K132").Value = 4695.081
$ws.Range("L132").Value = 11168.1819
$ws.Range("M132").Value = -2165.081
$ws.Range("N132").Value = -16228.1819

# row 136 (hunk @ diff line 28440)
$ws.Range("H136").Value = 9261668
$ws.Range("J136").Value = 27782106
$ws.Range("L136").Value = 83346318
$ws.Range("N136").Value = -83351418

$ws = $wb.Worksheets.Item("CUL")
# row 2 (hunk @ diff line 28831)
$ws.Range("H2").Value = 38.17647
$ws.Range("I2").Value = 21.666666
$ws.Range("J2").Value = 41.714287
$ws.Range("K2").Value = 129.999996
$ws.Range("L2").Value = 250.285722
$ws.Range("M2").Value = -16.99999600000001
$ws.Range("N2").Value = -476.285722

# row 87 (hunk @ diff line 33185)
$ws.Range("H87").Value = 12646.154
$ws.Range("J87").Value = 15677.777
$ws.Range("L87").Value = 47033.331
$ws.Range("N87").Value = -49529.331

# row 90 (hunk @ diff line 33338)
$ws.Range("H90").Value = 12646.154
$ws.Range("J90").Value = 15677.777
$ws.Range("L90").Value = 141099.993
$ws.Range("N90").Value = -153579.993

# row 120 (hunk @ diff line 34871)
$ws.Range("H120").Value = 20942.092

# row 124 (hunk @ diff line 35079)
$ws.Range("H124").Value = 14369.875
$ws.Range("I124").Value = 1952
$ws.Range("J124").Value = 35066.332
$ws.Range("K124").Value = 5856
$ws.Range("L124").Value = 105198.996
$ws.Range("M124").Value = -946
$ws.Range("N124").Value = -115018.996

# row 129 (hunk @ diff line 35333)
$ws.Range("H129").Value = 24299.13
$ws.Range("I129").Value = 2548.5715
$ws.Range("J129").Value = 58133.332
$ws.Range("K129").Value = 7645.7145
$ws.Range("L129").Value = 174399.996
$ws.Range("M129").Value = -2645.7145
$ws.Range("N129").Value = -184399.996

# row 131 (hunk @ diff line 35437)
$ws.Range("H131").Value = 1141.174
$ws.Range("I131").Value = 1304.6666
$ws.Range("J131").Value = 1083.4706
$ws.Range("K131").Value = 3913.9998
$ws.Range("L131").Value = 3250.4118
$ws.Range("M131").Value = 1126.0002
$ws.Range("N131").Value = -13330.4118

$ws = $wb.Worksheets.Item("LTW")
# row 136 (hunk @ diff line 49602)
$ws.Range("H136").Value = 3480.3044
$ws.Range("I136").Value = 3177.625
$ws.Range("J136").Value = 4172.143
$ws.Range("K136").Value = 9532.875
$ws.Range("L136").Value = 12516.429
$ws.Range("M136").Value = -6982.875
$ws.Range("N136").Value = -17616.429

$ws = $wb.Worksheets.Item("WVR")
# row 68 (hunk @ diff line 53239)
$ws.Range("H68").Value = 39654.2
$ws.Range("J68").Value = 39654.2
$ws.Range("L68").Value = 39654.2
$ws.Range("N68").Value = -41276.2

# row 71 (hunk @ diff line 53383)
$ws.Range("H71").Value = 39654.2
$ws.Range("J71").Value = 39654.2
$ws.Range("L71").Value = 118962.6
$ws.Range("N71").Value = -127074.6

# row 81 (hunk @ diff line 53870)
$ws.Range("H81").Value = 1470
$ws.Range("I81").Value = 1801
$ws.Range("J81").Value = 1139
$ws.Range("K81").Value = 3602
$ws.Range("L81").Value = 2278
$ws.Range("M81").Value = -2541
$ws.Range("N81").Value = -4400

# row 84 (hunk @ diff line 54017)
$ws.Range("H84").Value = 1470
$ws.Range("I84").Value = 1801
$ws.Range("J84").Value = 1139
$ws.Range("K84").Value = 18010
$ws.Range("L84").Value = 11390
$ws.Range("M84").Value = -12706
$ws.Range("N84").Value = -21998

# row 92 (hunk @ diff line 54397)
$ws.Range("H92").Value = 38950
$ws.Range("J92").Value = 38950
$ws.Range("L92").Value = 38950
$ws.Range("N92").Value = -43942

# row 107 (hunk @ diff line 55132)
$ws.Range("H107").Value = 1836.25
$ws.Range("I107").Value = 647.125
$ws.Range("K107").Value = 1941.375
$ws.Range("M107").Value = -21.375

# row 132 (hunk @ diff line 56360)
$ws.Range("H132").Value = 9393.343999999999
$ws.Range("I132").Value = 935.44684
$ws.Range("J132").Value = 32776.94
$ws.Range("K132").Value = 2806.34052
$ws.Range("L132").Value = 98330.82000000001
$ws.Range("M132").Value = -276.3405199999997
$ws.Range("N132").Value = -103390.82
